# [Report1] Minor stuff, cahier charges
#
# 1. The stray "_GoBack" bookmark (Word's "last edit position" marker)
#    was sitting at the very top of the cover page. It is removed from
#    there and re-created at the start of the "Pierre-Marie AIRIAU"
#    paragraph (the paragraph right after "Etudiants"), i.e. where the
#    author's cursor actually was the last time the document was saved.
# 2. The run holding "Etudiants" is split in two ("E" / "tudiants")
#    with identical run formatting - a harmless side effect of nudging
#    the bold formatting of the first letter.

$d = $word.ActiveDocument

# --- 1a. Drop the old _GoBack bookmark (top of the cover page) -------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Locate "Etudiants" and split its run into "E" + "tudiants" ---
$hit = $d.Content
$found = $hit.Find.Execute("Etudiants", $false, $false, $false, $false,
                            $false, $true, 1, $false, "", 0)

if ($found) {
    # Paragraph that contains "Etudiants" - needed to locate the
    # paragraph that immediately follows it.
    $etudiantsPara = $hit.Paragraphs(1)

    # Re-apply Bold on just the leading "E". Toggling it off then back
    # on keeps the visible formatting identical but forces the run to
    # split into two runs ("E" and "tudiants") that share the same rPr.
    $firstLetter = $d.Range($hit.Start, $hit.Start + 1)
    $firstLetter.Font.Bold = 0
    $firstLetter.Font.Bold = 1

    # --- 1b. Re-create _GoBack at the start of the next paragraph ----
    $nextPara = $etudiantsPara.Next()
    if ($nextPara -ne $null) {
        $target = $nextPara.Range
        $target.Collapse(1)
        $d.Bookmarks.Add("_GoBack", $target)
    }
}
